# Locate the very end of the document (right before the trailing paragraph
# mark of the last paragraph, which carries the "_GoBack" bookmark).
$d = $word.ActiveDocument
$endPos = $d.Content.End

# 1) Append the new "。。。。" text to the end of the document (it lands in
#    the same run/paragraph as the existing "，，，，，，" text, right before
#    the bookmark-carrying paragraph mark).
$insertPoint = $d.Range($endPos - 1, $endPos - 1)
$insertPoint.InsertBefore("。。。。")

# 2) Split the paragraph so "，，，，，，" and "。。。。" become separate
#    paragraphs (the bookmark stays with the paragraph mark, i.e. now with
#    the "。。。。" paragraph).
$splitPoint = $d.Range($endPos - 1, $endPos - 1)
$splitPoint.InsertBefore("`r")

# 3) Split again at the same spot to add a blank paragraph between
#    "，，，，，，" and "。。。。".
$splitPoint2 = $d.Range($endPos - 1, $endPos - 1)
$splitPoint2.InsertBefore("`r")
